# Auto-generated edit script for cryptos.xlsx update
# Updates Price (D) and Volume(1h) (E) columns for rows 2-51

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force Text format on the affected range so that numeric-looking
# strings (e.g. "246.19", "1.000") are stored as text, matching the source data,
# instead of being auto-converted to numbers by Excel.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = "30.594.94"
$ws.Range("D3").Value = "1.883.79"
$ws.Range("E3").Value = "  -0.59%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "246.19"
$ws.Range("E5").Value = "  -0.37%  "
$ws.Range("E6").Value = "  +0.08%  "
$ws.Range("D7").Value = "0.4741"
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("E8").Value = "  -1.24%  "
$ws.Range("D9").Value = "0.06540"
$ws.Range("E9").Value = "  +0.32%  "
$ws.Range("D10").Value = "22.35"
$ws.Range("E10").Value = "  -1.34%  "
$ws.Range("E11").Value = "  +2.80%  "
$ws.Range("D12").Value = "99.65"
$ws.Range("E12").Value = "  +2.76%  "
$ws.Range("D13").Value = "0.07824"
$ws.Range("E13").Value = "  +0.41%  "
$ws.Range("D14").Value = "1.881.84"
$ws.Range("E14").Value = "  -0.58%  "
$ws.Range("D15").Value = "5.234"
$ws.Range("E15").Value = "  -0.25%  "
$ws.Range("D16").Value = "283.77"
$ws.Range("E16").Value = "  -0.80%  "
$ws.Range("D17").Value = "30.569.61"
$ws.Range("E17").Value = "  -0.60%  "
$ws.Range("E18").Value = "  -0.80%  "
$ws.Range("D19").Value = "0.000007523"
$ws.Range("E19").Value = "  -0.16%  "
$ws.Range("D20").Value = "0.9998"
$ws.Range("D21").Value = "2.127.67"
$ws.Range("E21").Value = "  +0.09%  "
$ws.Range("D22").Value = "5.351"
$ws.Range("E22").Value = "  +0.24%  "
$ws.Range("E23").Value = "  +0.06%  "
$ws.Range("D24").Value = "6.434"
$ws.Range("E24").Value = "  +2.28%  "
$ws.Range("D25").Value = "9.170"
$ws.Range("E25").Value = "  -0.65%  "
$ws.Range("D26").Value = "163.85"
$ws.Range("E26").Value = "  -0.64%  "
$ws.Range("E27").Value = "  -0.14%  "
$ws.Range("E28").Value = "  -1.09%  "
$ws.Range("D29").Value = "0.09750"
$ws.Range("E29").Value = "  -0.56%  "
$ws.Range("E30").Value = "  -1.17%  "
$ws.Range("D31").Value = "1.503"
$ws.Range("E31").Value = "  +0.79%  "
$ws.Range("D32").Value = "4.249"
$ws.Range("E32").Value = "  -1.46%  "
$ws.Range("D33").Value = "4.185"
$ws.Range("E33").Value = "  -0.08%  "
$ws.Range("D34").Value = "0.04842"
$ws.Range("E34").Value = "  -1.37%  "
$ws.Range("D35").Value = "1.132"
$ws.Range("E35").Value = "  -0.28%  "
$ws.Range("D36").Value = "0.6983"
$ws.Range("E36").Value = "  -0.24%  "
$ws.Range("E37").Value = "  +2.35%  "
$ws.Range("D38").Value = "0.01904"
$ws.Range("E38").Value = "  +0.11%  "
$ws.Range("D39").Value = "2.871"
$ws.Range("E39").Value = "  +1.21%  "
$ws.Range("D40").Value = "6.305"
$ws.Range("E40").Value = "  -0.32%  "
$ws.Range("D41").Value = "75.49"
$ws.Range("E41").Value = "  -0.85%  "
$ws.Range("D42").Value = "1.975"
$ws.Range("E42").Value = "  -2.05%  "
$ws.Range("D43").Value = "0.4250"
$ws.Range("E43").Value = "  -1.30%  "
$ws.Range("D44").Value = "1.000"
$ws.Range("E44").Value = "  +0.03%  "
$ws.Range("D45").Value = "0.8383"
$ws.Range("E45").Value = "  +0.13%  "
$ws.Range("D46").Value = "9.962"
$ws.Range("E46").Value = "  +3.41%  "
$ws.Range("D47").Value = "101.41"
$ws.Range("E47").Value = "  -0.46%  "
$ws.Range("E48").Value = "  -0.40%  "
$ws.Range("D49").Value = "35.27"
$ws.Range("E49").Value = "  -0.51%  "
$ws.Range("D50").Value = "0.05775"
$ws.Range("E50").Value = "  +0.13%  "
$ws.Range("E51").Value = "  -0.61%  "

# Restore original (default) style so no stray per-cell formatting is introduced.
$dataRange.Style = "Normal"

Write-Host "Updated cryptos list"
